$d = $word.ActiveDocument

# The existing "_GoBack" bookmark (left over from the last edit location)
# currently sits at the end of the 3rd paragraph. Word will re-create it
# at the new last edit location once we type the new paragraph below, so
# remove the stale one first to avoid a duplicate-name collision.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# Append a new paragraph after the last paragraph in the document.
$lastPara = $d.Paragraphs.Last
$tailRange = $lastPara.Range
$tailRange.InsertParagraphAfter()

# Build the new paragraph's content as raw WordprocessingML so that the
# "github" token is wrapped the way Word's "check spelling as you type"
# would mark it (<w:proofErr w:type="spellStart"/> ... spellEnd), and so
# the _GoBack bookmark ends up collapsed at the very end of the new
# paragraph (the new "last edited" location), matching real Word's
# behaviour when you type at the end of a document.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:r><w:t xml:space="preserve">This is line 4 in </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>github</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t>/b1 branch</w:t></w:r>' +
       '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
       '<w:bookmarkEnd w:id="0"/>' +
       '</w:p>'
$newRange.InsertXML($xml)
